$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook previously had two blocks of three rows each: one where the
# sending cluster was "ECs" (rows 2-4) and one where it was "MuSCs"
# (rows 5-7). The TPM recomputation drops the "ECs" sender entirely, so
# those three rows are removed and the "MuSCs" block shifts up to become
# the only data rows (2-4), with refreshed TPM-derived figures.
$ws.Rows.Item(2).Resize(3).Delete() | Out-Null

# Row 2: MuSCs -> Slitrk2 -> Ptprs -> ECs
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1841803333333333
$ws.Range("H2").Value = 0.5525409999999999
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.425633666666667
$ws.Range("N2").Value = 7.276901000000001
$ws.Range("O2").Value = 0.0662600404061536
$ws.Range("P2").Value = 0.06626004040615362
$ws.Range("Q2").Value = 0.4467540172712222
$ws.Range("R2").Value = 4.020786155441
$ws.Range("S2").Value = 0.0662600404061536
$ws.Range("T2").Value = 0.06626004040615362

# Row 3: MuSCs -> Slitrk2 -> Ptprs -> FAPs
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1841803333333333
$ws.Range("H3").Value = 0.5525409999999999
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("O3").Value = 0.4234968256437875
$ws.Range("P3").Value = 0.4234968256437876
$ws.Range("Q3").Value = 2.855399830701
$ws.Range("R3").Value = 25.698598476309
$ws.Range("S3").Value = 0.4234968256437875
$ws.Range("T3").Value = 0.4234968256437876

# Row 4: MuSCs -> Slitrk2 -> Ptprs -> MuSCs
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1841803333333333
$ws.Range("H4").Value = 0.5525409999999999
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 18.67887366666667
$ws.Range("N4").Value = 56.03662100000001
$ws.Range("O4").Value = 0.5102431339500588
$ws.Range("P4").Value = 0.5102431339500588
$ws.Range("Q4").Value = 3.440281178217889
$ws.Range("R4").Value = 30.962530603961
$ws.Range("S4").Value = 0.5102431339500588
$ws.Range("T4").Value = 0.5102431339500588

Write-Host "Applied TPM update: removed ECs-sender rows, refreshed MuSCs-sender rows."
